$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C: "Hoodtape Volume 2" -> "Hoodtape volume 2" (rows 2-29) ---
# --- 2) Column F: update credits metadata block (rows 2-29) ---

$oldTail = "Producers`n`n    `n      Alexis Troy, B-Case, Beatzarre & 10 more`nWriters`n`n    `n      David Ruoff, Elias Klughammer & Kollegah`nCopyright ©`n`n    `n      Alpha Music Empire`nLabel`n`n    `n      Alpha Music Empire`nPhonographic Copyright ℗`n`n    `n      Alpha Music Empire"
$newTail = "Producers`n`n    `n      Alexis Troy, B-Case, Beatzarre & 11 more`nWriters`n`n    `n      Alexis Troy, B-Case, Beatzarre & 12 more`nComposer`n`n    `n      Alexis Troy, B-Case, Beatzarre & 11 more`nCopyright ©`n`n    `n      Alpha Music Empire`nLabel`n`n    `n      Alpha Music Empire"

for ($r = 2; $r -le 29; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cText = $cCell.Text
    if ($cText -eq "Hoodtape Volume 2") {
        $cCell.Value = "Hoodtape volume 2"
    }

    $fCell = $ws.Cells.Item($r, 6)
    $fText = $fCell.Text
    if ($fText.Contains($oldTail)) {
        $fCell.Value = $fText.Replace($oldTail, $newTail)
    }
}

# --- 3) Row 28, column D: merge two lyric lines into one ---
$dCell = $ws.Cells.Item(28, 4)
$dText = $dCell.Text
$oldLines = "Man tauscht Gegenstände wie Eheringe und Goldketten`nFür Medizin, wesentliche Lebensmittel und Wolldecken"
$newLines = "Man tauscht Gegenstände wie Eheringe und GoldkettenFür Medizin, wesentliche Lebensmittel und Wolldecken"
if ($dText.Contains($oldLines)) {
    $dCell.Value = $dText.Replace($oldLines, $newLines)
}
